$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E ("N° de lot") only ever contained a blank/whitespace placeholder
# string for the data rows. Clear that column's contents for rows 2-43 so the
# now-unused shared string gets dropped from the workbook on save.
$ws.Range("E2:E43").ClearContents()
